$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '35.808.88'
$ws.Range("E2").Value = '  +3.87%  '

$ws.Range("D3").Value = '1.873.42'
$ws.Range("E3").Value = '  +3.33%  '

$ws.Range("D5").Value = '''233.08'
$ws.Range("E5").Value = '  +3.33%  '

$ws.Range("E6").Value = '  +3.31%  '

$ws.Range("D7").Value = '''1.01'
$ws.Range("E7").Value = '  +0.42%  '

$ws.Range("D8").Value = '''42.39'
$ws.Range("E8").Value = '  +10.57%  '

$ws.Range("E9").Value = '  +7.72%  '

$ws.Range("E10").Value = '  +3.66%  '

$ws.Range("E11").Value = '  +4.09%  '

$ws.Range("D12").Value = '2.146.91'
$ws.Range("E12").Value = '  +3.47%  '

$ws.Range("D13").Value = '''11.69'
$ws.Range("E13").Value = '  +4.42%  '

$ws.Range("B14").Value = 'Polygon'
$ws.Range("C14").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D14").Value = '''0.686'
$ws.Range("E14").Value = '  +8.44%  '

$ws.Range("B15").Value = 'WrappedEther'
$ws.Range("C15").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D15").Value = '1.863.79'
$ws.Range("E15").Value = '  +2.71%  '

$ws.Range("E16").Value = '  +8.42%  '

$ws.Range("D17").Value = '35.844.73'
$ws.Range("E17").Value = '  +3.99%  '

$ws.Range("D18").Value = '''70.84'
$ws.Range("E18").Value = '  +3.79%  '

$ws.Range("D19").Value = '0.0₃0811'
$ws.Range("E19").Value = '  +4.80%  '

$ws.Range("D20").Value = '''249.00'
$ws.Range("E20").Value = '  +2.65%  '

$ws.Range("D21").Value = '''12.50'
$ws.Range("E21").Value = '  +11.30%  '

$ws.Range("D22").Value = '''4.83'
$ws.Range("E22").Value = '  +16.97%  '

$ws.Range("D24").Value = '''2.25'
$ws.Range("E24").Value = '  +1.79%  '

$ws.Range("D25").Value = '''171.81'
$ws.Range("E25").Value = '  +0.95%  '

$ws.Range("D26").Value = '''8.09'
$ws.Range("E26").Value = '  +3.46%  '

$ws.Range("D27").Value = '''18.03'
$ws.Range("E27").Value = '  +2.78%  '

$ws.Range("D28").Value = '''0.124'
$ws.Range("E28").Value = '  +2.37%  '

$ws.Range("E29").Value = '  +17.64%  '

$ws.Range("E30").Value = '  +0.37%  '

$ws.Range("D31").Value = '3.370.87'
$ws.Range("E31").Value = '  +38.74%  '

$ws.Range("D32").Value = '''0.0554'
$ws.Range("E32").Value = '  +7.13%  '

$ws.Range("E33").Value = '  +4.98%  '

$ws.Range("E34").Value = '  +6.79%  '

$ws.Range("D35").Value = '''1.92'
$ws.Range("E35").Value = '  +5.44%  '

$ws.Range("D36").Value = '''98.58'
$ws.Range("E36").Value = '  +20.60%  '

$ws.Range("E37").Value = '  +7.39%  '

$ws.Range("E38").Value = '  +7.58%  '

$ws.Range("D39").Value = '1.366.23'
$ws.Range("E39").Value = '  +0.26%  '

$ws.Range("E40").Value = '  +3.20%  '

$ws.Range("E41").Value = '  +6.00%  '

$ws.Range("E42").Value = '  +8.26%  '

$ws.Range("D43").Value = '''15.16'
$ws.Range("E43").Value = '  +10.03%  '

$ws.Range("E44").Value = '  +3.28%  '

$ws.Range("E45").Value = '  +1.97%  '

$ws.Range("E46").Value = '  +1.24%  '

$ws.Range("D47").Value = '''6.37'
$ws.Range("E47").Value = '  +10.44%  '

$ws.Range("D48").Value = '''0.0521'
$ws.Range("E48").Value = '  +2.14%  '

$ws.Range("D49").Value = '2.043.78'
$ws.Range("E49").Value = '  +3.45%  '

$ws.Range("D50").Value = '''105.55'
$ws.Range("E50").Value = '  +3.48%  '

$ws.Range("E51").Value = '  +0.41%  '
